# "Sign up Request Meeting" sheet: drop the "url" column (column A) entirely.
# Excel shifts title/firstName/lastName/email/password/confirmPassword (and
# their row-2 values) one column to the left, the mailto hyperlink that used
# to live on E2 ends up on D2, and the old hyperlink to the sandbox URL
# (which lived on A2) goes away with the deleted column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the whole "url" column; everything to its right shifts left by one.
$ws.Columns.Item(1).Delete()

# Re-anchor the surviving "confirmPassword"/email hyperlink cleanly on its
# new home (D2) with the standard Hyperlink cell style, and drop the stale
# hyperlink metadata left over from the deleted url column.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:dennisgituto@yahoo.com")
$ws.Range("D2").Style = "Hyperlink"

# Reset the active selection back to the top-left cell.
$ws.Range("A1").Select() | Out-Null
